# "Let's continue selling Rocket engines for a while"
#
# The BFR launchpad level input (C27) is bumped on both the "Production
# Plans second stage" and "Production Plans final stage" sheets, and the
# Rocket engine market price (C4) is lowered on the final-stage sheet.
# Everything else (D3:F10, C18, I18:I20, C22:C23, I22:I25, F27, I27:I33,
# C35:C38/C36:C38, C44:C46/C45:C46, etc.) is formula-derived and
# recalculates automatically.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Input cell edits
# ---------------------------------------------------------------------

$wsSecond = $wb.Worksheets.Item("Production Plans second stage")
$wsSecond.Range("C27").Value = 0.31121951219512189

$wsFinal = $wb.Worksheets.Item("Production Plans final stage")
$wsFinal.Range("C4").Value = 900
$wsFinal.Range("C27").Value = 0.28292682926829271

# ---------------------------------------------------------------------
# 2. Sheet view / selection updates
# ---------------------------------------------------------------------

# "Rocket production" loses the tab-selected / scrolled state and the
# selection moves from D20 to E20.
$wsRocket = $wb.Worksheets.Item("Rocket production")
$wsRocket.Activate()
$excel.Goto($wsRocket.Range("E20"), $false)

# "Production Plans second stage": scrolled a bit further down/right and
# selection moves from F36 to D36.
$wsSecond.Activate()
$wsSecond.Range("D36").Select()

# "Production Plans final stage" becomes the active tab, scrolled back to
# the top-left area, with the selection moved from G19 to E42.
$wsFinal.Activate()
$excel.Goto($wsFinal.Range("E42"), $false)
